$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Session 21) is being filled in: date plan/actual both set to 2025.01.13,
# progress bumped to 90%, subject "SVM" (new) and problem "too many formulas" added.
# That also means the row is no longer "pending" so the yellow highlight on the
# Duration/Plan cells is cleared (matching the look of the already-completed rows above it).
$ws.Range("B8").Interior.Pattern = -4142
$ws.Range("C8").Interior.Pattern = -4142

$ws.Range("D8").Value = $ws.Range("C8").Value
$ws.Range("E8").Value = 0.9
$ws.Range("F8").Value = "SVM"
$ws.Range("G8").Value = "too many formulas"

# Move the active selection the way it ended up after the edit.
$ws.Range("E9").Select()
